$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Address") to hold the new
# "Company Name" field, shifting all existing columns (B..AE) right by one
# (B..AF).
$ws.Columns("B").EntireColumn.Insert() | Out-Null

# Populate the new column B (order chosen so the shared-string table is
# built up in the same sequence as the authored workbook).
$ws.Range("B2").Value = "{vendor:company_name}"

# Add the new trailing column AG for the vendor active/inactive status.
$ws.Range("AG1").Value = "Status"
$ws.Range("AG2").Value = "{vendor:active_status}"

$ws.Range("B1").Value = "Company Name"

# Match column widths to the new layout (column B mirrors column A, the
# new Status column gets its own width).
$ws.Range("B1").ColumnWidth = $ws.Range("A1").ColumnWidth
$ws.Range("AG1").ColumnWidth = 20.5

# Move the active selection as it appears in the saved workbook.
$ws.Range("B5").Select() | Out-Null
